$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new participant's data
$ws.Range("A2").Value = "1MM43"
$ws.Range("B2").Value = "Michael Jackson"
$ws.Range("C2").Value = 43
$ws.Range("D2").Value = "Masculino"
$ws.Range("E2").Value = "30/12/2020 - 4:13:30 p. m."
$ws.Range("F2").Value = 5558545

# Remove the other two participant rows (rows 3 and 4)
$ws.Range("A3:F4").Delete()

# Update selection to match the author's final cursor/selection state
$ws.Range("A2:F10").Select()
